$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2's value (C2 formula will recalc automatically)
$ws.Range("B2").Value = 7

# Update the active selection to G4
$ws.Range("G4").Select()
